# "contingencies with rene fine"
# Update line/extraction results for rows 8-15 and append two new lines
# (line7, line8) as rows 16-17.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 8 (line index 6) ---
$ws.Range("C8").Value = 14
$ws.Range("D8").Value = 11
$ws.Range("E8").Value = $true

# --- Row 9 (line index 7) ---
$ws.Range("C9").Value = 16
$ws.Range("E9").Value = $true

# --- Row 10 (line index 8) ---
$ws.Range("C10").Value = 5
$ws.Range("D10").Value = 12
$ws.Range("E10").Value = $true

# --- Row 11 (line index 9) ---
$ws.Range("C11").Value = 5
$ws.Range("D11").Value = 9

# --- Row 12 (line index 10) ---
$ws.Range("C12").Value = 10

# --- Row 13 (line index 11) ---
$ws.Range("D13").Value = 8

# --- Row 14 (line index 12) ---
$ws.Range("C14").Value = 9
$ws.Range("D14").Value = 11
$ws.Range("E14").Value = $true

# --- Row 15 (line index 13) ---
$ws.Range("C15").Value = 7
$ws.Range("D15").Value = 11
$ws.Range("E15").Value = $true

# --- Row 16 (new, line index 14 / "line7") ---
# Clone the formatting of A15 (bold + bordered style) onto A16 before
# writing its value, so the new index cell keeps the same cell style.
$ws.Range("A15").Copy()
$ws.Range("A16").PasteSpecial(-4122)
$ws.Range("A16").Value = 14
$ws.Range("B16").Value = "line7"
$ws.Range("C16").Value = 5
$ws.Range("D16").Value = 7
$ws.Range("E16").Value = $true

# --- Row 17 (new, line index 15 / "line8") ---
$ws.Range("A15").Copy()
$ws.Range("A17").PasteSpecial(-4122)
$ws.Range("A17").Value = 15
$ws.Range("B17").Value = "line8"
$ws.Range("C17").Value = 8
$ws.Range("D17").Value = 5
$ws.Range("E17").Value = $true

$excel.CutCopyMode = $false
